$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the "Date" paragraph timestamp.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "June   8, 2021 (09:17:30 PM)", $true, $false, $false, $false, $false,
    $true, 1, $false, "June   8, 2021 (10:31:14 PM)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new bulleted hyperlink paragraph
#    "https://www.w3schools.com/cs/trycs.asp?filename=demo_helloworld"
#    right after the "https://dotnetfiddle.net/" bullet and before the
#    "Note that none of them are endorsed..." paragraph, in the
#    "Compiling Code On-Line" section.
# ---------------------------------------------------------------------------
$dotnetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*dotnetfiddle.net*") {
        $dotnetPara = $p
    }
}

$endOfDotnet = $dotnetPara.Range.End
$insertionPoint = $d.Range($endOfDotnet, $endOfDotnet)
$insertionPoint.InsertAfter("w3schools-placeholder`r")

# Find the freshly created (still empty-ish) paragraph that now sits between
# the dotnetfiddle bullet and the "Note that..." paragraph.
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*w3schools-placeholder*") {
        $newPara = $p
    }
}

# Give it the same "Compact" style + bullet numbering (numId 1006) as its
# sibling links.
$newPara.Style = "Compact"
$listTemplate = $dotnetPara.Range.ListFormat.ListTemplate
$newPara.Range.ListFormat.ApplyListTemplate($listTemplate, $true)

# Turn its text into the real hyperlink (drop the trailing paragraph mark
# from the target range so the hyperlink wraps only the visible text).
$newPara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*w3schools-placeholder*") {
        $newPara2 = $p
    }
}
$hyperlinkRange = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
$d.Hyperlinks.Add(
    $hyperlinkRange,
    "https://www.w3schools.com/cs/trycs.asp?filename=demo_helloworld",
    "",
    "",
    "https://www.w3schools.com/cs/trycs.asp?filename=demo_helloworld") | Out-Null
